$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose trials fell out of the refreshed query result
# (row 7 = CADANCE, row 5 = KATALYST). Delete the lower one first so the
# higher row index still points at the right row.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()

# Update values changed by the refreshed query result
$ws.Cells.Item(2, 2).Value = 50
$ws.Cells.Item(13, 2).Value = 0

# Clear the leftover numeric-format style from the trial-name column so it
# matches the plain "General" style the refresh produced
$ws.Range("A2:A13").Style = "Normal"

# Keep the hidden ExternalData_1 defined name (used by the query table) in
# sync with the now-smaller data range
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$B`$13"
